# Rename DPV version namespace from "v2.0" to "2.0" (drop leading "v")
# in the "Namespaces-v2" sheet's namespace column (column B, rows 2-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Namespaces-v2")

$lastRow = 25
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $value = $cell.Value2
    if ($value -ne $null -and $value -like "*v2.0*") {
        $cell.Value2 = $value -replace "v2\.0", "2.0"
    }
}
